$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Add the new log entry as row 19
$ws.Range("A19").Value = "Kan mijn wachtwoord niet resetten"
$ws.Range("B19").Value = "mailmind.test@zohomail.eu"
$ws.Range("C19").Value = "Ik krijg geen e-mail bij wachtwoord resetten."
$ws.Range("D19").Value = "IT / Technisch probleem"
$ws.Range("E19").Value = "Beste klant,`nBedankt voor uw bericht. Het spijt me te horen dat u geen e-mail heeft ontvangen bij het resetten van uw wachtwoord. Om dit probleem zo snel mogelijk voor u op te lossen, zou u ons alstublieft uw gebruikersnaam willen doorgeven? Op die manier kunnen we gerichter onderzoeken waar het probleem zich voordoet en u verder helpen.`nIk kijk ernaar uit om uw gegevens te ontvangen, zodat we dit snel kunnen oplossen.`nMet vriendelijke groet,`n[Naam] E-mailassistent at [Bedrijfsnaam]"
$ws.Range("F19").Value = "2025-06-22 18:49:12"
$ws.Range("G19").Value = "Ja"

# Undo the automatic row-height expansion triggered by the multi-line
# text in E19 so the row matches the sheet's default (no explicit height).
$ws.Rows.Item(19).AutoFit()

# Extend conditional formatting ranges to include the new row
$cfD = $ws.Range("D2:D18").FormatConditions
for ($i = 1; $i -le $cfD.Count(); $i++) {
    $cfD.Item($i).ModifyAppliesToRange($ws.Range("D2:D19"))
}

$cfG = $ws.Range("G2:G18").FormatConditions
for ($i = 1; $i -le $cfG.Count(); $i++) {
    $cfG.Item($i).ModifyAppliesToRange($ws.Range("G2:G19"))
}

# Update the Dashboard summary count for "IT / Technisch probleem" (3 -> 4)
$dash.Range("B2").Value = 4
